$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3980872333049774
$ws.Range("B1").Value = 3.237294912338257
$ws.Range("C1").Value = 4.636600494384766
$ws.Range("D1").Value = 1.840664982795715
$ws.Range("E1").Value = 0.8051331043243408
